$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 19:55"

# --- Update per-country statistics (columns B..H) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 8476461
$ws.Range("C4").Value = 19808
$ws.Range("D4").Value = 5515763
$ws.Range("E4").Value = 2735072
$ws.Range("G4").Value = 404
$ws.Range("H4").Value = 225626

# Row 5: India
$ws.Range("B5").Value = 7645741
$ws.Range("C5").Value = 51005
$ws.Range("D5").Value = 6788363
$ws.Range("E5").Value = 741490
$ws.Range("G5").Value = 652
$ws.Range("H5").Value = 115888

# Row 8: Espana
$ws.Range("B8").Value = 1029668
$ws.Range("C8").Value = 13873
$ws.Range("G8").Value = 218
$ws.Range("H8").Value = 34210

# Row 21: Alemania
$ws.Range("B21").Value = 377278
$ws.Range("C21").Value = 3547
$ws.Range("E21").Value = 69054
$ws.Range("G21").Value = 25
$ws.Range("H21").Value = 9924

# Row 24: Turquia
$ws.Range("B24").Value = 351413
$ws.Range("C24").Value = 1894
$ws.Range("D24").Value = 306939
$ws.Range("E24").Value = 35029
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 9445

# Row 28: Israel
$ws.Range("B28").Value = 305993
$ws.Range("C28").Value = 1117
$ws.Range("D28").Value = 280849
$ws.Range("E28").Value = 22866
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = 2278

# Row 35: Marruecos
$ws.Range("B35").Value = 179003
$ws.Range("C35").Value = 3254
$ws.Range("D35").Value = 148838
$ws.Range("E35").Value = 27138
$ws.Range("G35").Value = 51
$ws.Range("H35").Value = 3027

# Row 52: Etiopia
$ws.Range("B52").Value = 90490
$ws.Range("C52").Value = 630
$ws.Range("D52").Value = 43638
$ws.Range("E52").Value = 45481
$ws.Range("G52").Value = 6
$ws.Range("H52").Value = 1371

# Row 67: Argelia
$ws.Range("B67").Value = 54829
$ws.Range("C67").Value = 213
$ws.Range("D67").Value = 38346
$ws.Range("E67").Value = 14610
$ws.Range("G67").Value = 8
$ws.Range("H67").Value = 1873

# Rows 125/126: Sri Lanka overtakes Suazilandia (Eswatini) in the sort order.
# Suazilandia's figures stay exactly the same as before; Sri Lanka's figures
# are refreshed and now exceed Suazilandia's, so the two rows swap places.
$ws.Range("A125").Value = "Sri Lanka"
$ws.Range("B125").Value = 5811
$ws.Range("C125").Value = 186
$ws.Range("D125").Value = 3457
$ws.Range("E125").Value = 2341
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 13

$ws.Range("A126").Value = "Suazilandia"
$ws.Range("B126").Value = 5788
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 5427
$ws.Range("E126").Value = 245
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 116

# Row 181: Comoras
$ws.Range("B181").Value = 504
$ws.Range("C181").Value = 2
$ws.Range("D181").Value = 494
$ws.Range("E181").Value = 3
